$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10.406695138989
$ws.Range("D2").Value = 9.132519501172135
$ws.Range("E2").Value = 13.6846271886964
$ws.Range("F2").Value = 32.97028064821509
$ws.Range("G2").Value = 34.50059769325176
$ws.Range("H2").Value = 15.90999113674322
$ws.Range("I2").Value = 25.83392270719255
$ws.Range("J2").Value = 10.223115823806
$ws.Range("K2").Value = 15.26217302007433
$ws.Range("L2").Value = 10.37964113127259
$ws.Range("M2").Value = 17.73752886416322
$ws.Range("O2").Value = 24.94338156223582
$ws.Range("C3").Value = 10.37738794702658
$ws.Range("D3").Value = 9.1084069595149
$ws.Range("E3").Value = 13.70990106527195
$ws.Range("F3").Value = 33.12286436698575
$ws.Range("G3").Value = 34.72523185142013
$ws.Range("H3").Value = 15.98074198610993
$ws.Range("I3").Value = 25.95344731472784
$ws.Range("J3").Value = 10.24837593625777
$ws.Range("K3").Value = 14.70484470738192
$ws.Range("L3").Value = 10.40324589937164
$ws.Range("M3").Value = 17.50988231451848
$ws.Range("O3").Value = 25.07721638226649
$ws.Range("C4").Value = 10.36088196678621
$ws.Range("D4").Value = 9.094587685460294
$ws.Range("E4").Value = 13.72712394443978
$ws.Range("F4").Value = 33.22437503076451
$ws.Range("G4").Value = 34.87404028896091
$ws.Range("H4").Value = 16.02684501451412
$ws.Range("I4").Value = 26.03237175636675
$ws.Range("J4").Value = 10.26474265736157
$ws.Range("K4").Value = 14.35114582578135
$ws.Range("L4").Value = 10.41859956433082
$ws.Range("M4").Value = 17.36942496697583
$ws.Range("O4").Value = 25.1648637131977
$ws.Range("C5").Value = 10.35453478299798
$ws.Range("D5").Value = 9.089207677334786
$ws.Range("E5").Value = 13.73457147284757
$ws.Range("F5").Value = 33.26770514606132
$ws.Range("G5").Value = 34.9374079089024
$ws.Range("H5").Value = 16.04630231622122
$ws.Range("I5").Value = 26.06592478865971
$ws.Range("J5").Value = 10.27162825979689
$ws.Range("K5").Value = 14.2042997563945
$ws.Range("L5").Value = 10.42507313611618
$ws.Range("M5").Value = 17.31206759493847
$ws.Range("O5").Value = 25.20195607985992
$ws.Range("C6").Value = 10.35350387029873
$ws.Range("D6").Value = 9.08832961843661
$ws.Range("E6").Value = 13.73583405579673
$ws.Range("F6").Value = 33.27501857501833
$ws.Range("G6").Value = 34.94809448290076
$ws.Range("H6").Value = 16.04957366824407
$ws.Range("I6").Value = 26.07158019882243
$ws.Range("J6").Value = 10.2727846728677
$ws.Range("K6").Value = 14.17975800017957
$ws.Range("L6").Value = 10.42616117883689
$ws.Range("M6").Value = 17.30253771989876
$ws.Range("O6").Value = 25.2081982897127
$ws.Range("C7").Value = 10.3607948252503
$ws.Range("D7").Value = 9.094514106095795
$ws.Range("E7").Value = 13.72722264659112
$ws.Range("F7").Value = 33.22495144972434
$ws.Range("G7").Value = 34.87488385782866
$ws.Range("H7").Value = 16.02710470906541
$ws.Range("I7").Value = 26.03281863525137
$ws.Range("J7").Value = 10.26483464356716
$ws.Range("K7").Value = 14.34917613213115
$ws.Range("L7").Value = 10.4186859905684
$ws.Range("M7").Value = 17.36865184241396
$ws.Range("O7").Value = 25.16535838585552
$ws.Range("C8").Value = 10.39628398665708
$ws.Range("D8").Value = 9.124003042531086
$ws.Range("E8").Value = 13.69298808979918
$ws.Range("F8").Value = 33.02126554795758
$ws.Range("G8").Value = 34.57578726352455
$ws.Range("H8").Value = 15.93383396844764
$ws.Range("I8").Value = 25.87398501215189
$ws.Range("J8").Value = 10.23164799286823
$ws.Range("K8").Value = 15.07249442918854
$ws.Range("L8").Value = 10.38760182164389
$ws.Range("M8").Value = 17.65920851413289
$ws.Range("O8").Value = 24.98839164328457
$ws.Range("C9").Value = 10.47746803527749
$ws.Range("D9").Value = 9.189491605820638
$ws.Range("E9").Value = 13.63936045775201
$ws.Range("F9").Value = 32.68407598978057
$ws.Range("G9").Value = 34.07602962226095
$ws.Range("H9").Value = 15.77202097362168
$ws.Range("I9").Value = 25.60649759054673
$ws.Range("J9").Value = 10.17334294645611
$ws.Range("K9").Value = 16.39275926348912
$ws.Range("L9").Value = 10.33344814749955
$ws.Range("M9").Value = 18.22135072032041
$ws.Range("O9").Value = 24.6848098638044
$ws.Range("C10").Value = 10.5438630888342
$ws.Range("D10").Value = 9.242049989153998
$ws.Range("E10").Value = 13.6081675615132
$ws.Range("F10").Value = 32.47450284382711
$ws.Range("G10").Value = 33.76234117259202
$ws.Range("H10").Value = 15.66595387073108
$ws.Range("I10").Value = 25.43687137430485
$ws.Range("J10").Value = 10.13460024885442
$ws.Range("K10").Value = 17.29490055476328
$ws.Range("L10").Value = 10.29777576560381
$ws.Range("M10").Value = 18.62677137664884
$ws.Range("O10").Value = 24.48829784103403
$ws.Range("C11").Value = 10.57546008230016
$ws.Range("D11").Value = 9.266873560106363
$ws.Range("E11").Value = 13.59575317163278
$ws.Range("F11").Value = 32.38749524122212
$ws.Range("G11").Value = 33.63137800426414
$ws.Range("H11").Value = 15.62047610916824
$ws.Range("I11").Value = 25.36556271903203
$ws.Range("J11").Value = 10.11785663550453
$ws.Range("K11").Value = 17.68920979929163
$ws.Range("L11").Value = 10.28243397820069
$ws.Range("M11").Value = 18.80896156733471
$ws.Range("O11").Value = 24.40467046668352
$ws.Range("C12").Value = 10.58761878465256
$ws.Range("D12").Value = 9.276400467792401
$ws.Range("E12").Value = 13.59130692097155
$ws.Range("F12").Value = 32.35574925273362
$ws.Range("G12").Value = 33.58348416468881
$ws.Range("H12").Value = 15.60365306616083
$ws.Range("I12").Value = 25.33940363908048
$ws.Range("J12").Value = 10.11164233567956
$ws.Range("K12").Value = 17.83611791126552
$ws.Range("L12").Value = 10.27675129624837
$ws.Range("M12").Value = 18.87758215716698
$ws.Range("O12").Value = 24.37383340426744
$ws.Range("C13").Value = 10.58499169048931
$ws.Range("D13").Value = 9.274343117118299
$ws.Range("E13").Value = 13.59225317606351
$ws.Range("F13").Value = 32.36253279624516
$ws.Range("G13").Value = 33.59372320831945
$ws.Range("H13").Value = 15.60725849525925
$ws.Range("I13").Value = 25.34499989762184
$ws.Range("J13").Value = 10.11297509384652
$ws.Range("K13").Value = 17.8045871202403
$ws.Range("L13").Value = 10.27796952581204
$ws.Range("M13").Value = 18.86282074440441
$ws.Range("O13").Value = 24.38043774171515
$ws.Range("C14").Value = 10.57645654563012
$ws.Range("D14").Value = 9.267654831750615
$ws.Range("E14").Value = 13.59538227208662
$ws.Range("F14").Value = 32.38485936719024
$ws.Range("G14").Value = 33.62740363082326
$ws.Range("H14").Value = 15.61908408448282
$ws.Range("I14").Value = 25.3633936700507
$ws.Range("J14").Value = 10.11734285601463
$ws.Range("K14").Value = 17.70134478699061
$ws.Range("L14").Value = 10.2819639194367
$ws.Range("M14").Value = 18.81461473288037
$ws.Range("O14").Value = 24.40211682329274
$ws.Range("C15").Value = 10.57125353068642
$ws.Range("D15").Value = 9.263574428527889
$ws.Range("E15").Value = 13.59733210230597
$ws.Range("F15").Value = 32.39869168692859
$ws.Range("G15").Value = 33.64825547582087
$ws.Range("H15").Value = 15.62637947379829
$ws.Range("I15").Value = 25.37477035491466
$ws.Range("J15").Value = 10.12003464954048
$ws.Range("K15").Value = 17.63778962382442
$ws.Range("L15").Value = 10.28442711629504
$ws.Range("M15").Value = 18.78503743035561
$ws.Range("O15").Value = 24.41550412941668
$ws.Range("C16").Value = 10.54182561859785
$ws.Range("D16").Value = 9.240445710262257
$ws.Range("E16").Value = 13.60901455550107
$ws.Range("F16").Value = 32.48035695037353
$ws.Range("G16").Value = 33.77113705241071
$ws.Range("H16").Value = 15.66898171579316
$ws.Range("I16").Value = 25.44164955345302
$ws.Range("J16").Value = 10.13571216188813
$ws.Range("K16").Value = 17.26879944964527
$ws.Range("L16").Value = 10.29879617254692
$ws.Range("M16").Value = 18.61481555594639
$ws.Range("O16").Value = 24.49387926116001
$ws.Range("C17").Value = 10.52412472285614
$ws.Range("D17").Value = 9.226487886973901
$ws.Range("E17").Value = 13.61663574705365
$ws.Range("F17").Value = 32.53259208994608
$ws.Range("G17").Value = 33.84953479015368
$ws.Range("H17").Value = 15.69582673402561
$ws.Range("I17").Value = 25.48417883480721
$ws.Range("J17").Value = 10.14555500848078
$ws.Range("K17").Value = 17.03824506115709
$ws.Range("L17").Value = 10.30783766460344
$ws.Range("M17").Value = 18.50978136733073
$ws.Range("O17").Value = 24.54343798029289
$ws.Range("C18").Value = 10.5140752258814
$ws.Range("D18").Value = 9.218546086496787
$ws.Range("E18").Value = 13.62118639669634
$ws.Range("F18").Value = 32.56342000150988
$ws.Range("G18").Value = 33.89573097752569
$ws.Range("H18").Value = 15.71152827451195
$ws.Range("I18").Value = 25.50919154719647
$ws.Range("J18").Value = 10.15129927197849
$ws.Range("K18").Value = 16.90412766742323
$ws.Range("L18").Value = 10.31312149008362
$ws.Range("M18").Value = 18.44916069504458
$ws.Range("O18").Value = 24.57248555902748
$ws.Range("C19").Value = 10.510695442179
$ws.Range("D19").Value = 9.215872105651107
$ws.Range("E19").Value = 13.62275588978419
$ws.Range("F19").Value = 32.57399226197161
$ws.Range("G19").Value = 33.91156145551201
$ws.Range("H19").Value = 15.71688938264348
$ws.Range("I19").Value = 25.51775499530488
$ws.Range("J19").Value = 10.15325843727223
$ws.Range("K19").Value = 16.8584618608826
$ws.Range("L19").Value = 10.3149248401784
$ws.Range("M19").Value = 18.42860141271501
$ws.Range("O19").Value = 24.58241373640592
$ws.Range("C20").Value = 10.52599544270221
$ws.Range("D20").Value = 9.227964815999053
$ws.Range("E20").Value = 13.61580716313124
$ws.Range("F20").Value = 32.52695043495993
$ws.Range("G20").Value = 33.84107488612891
$ws.Range("H20").Value = 15.69294202624251
$ws.Range("I20").Value = 25.47959447994324
$ws.Range("J20").Value = 10.14449864247895
$ws.Range("K20").Value = 17.06294484184466
$ws.Range("L20").Value = 10.30686655420065
$ws.Range("M20").Value = 18.52098430190999
$ws.Range("O20").Value = 24.53810619028704
$ws.Range("C21").Value = 10.57895832543121
$ws.Range("D21").Value = 9.269615941566798
$ws.Range("E21").Value = 13.59445626936004
$ws.Range("F21").Value = 32.37826885867137
$ws.Range("G21").Value = 33.61746466945191
$ws.Range("H21").Value = 15.61559981219636
$ws.Range("I21").Value = 25.35796805236314
$ws.Range("J21").Value = 10.11605651878738
$ws.Range("K21").Value = 17.73173562318274
$ws.Range("L21").Value = 10.28078722820145
$ws.Range("M21").Value = 18.82878445563497
$ws.Range("O21").Value = 24.39572659100377
$ws.Range("C22").Value = 10.61469780182196
$ws.Range("D22").Value = 9.297573916250313
$ws.Range("E22").Value = 13.58198719989612
$ws.Range("F22").Value = 32.2881044314409
$ws.Range("G22").Value = 33.48123137923421
$ws.Range("H22").Value = 15.56737425060125
$ws.Range("I22").Value = 25.28339833876916
$ws.Range("J22").Value = 10.09820298153849
$ws.Range("K22").Value = 18.1547562986336
$ws.Range("L22").Value = 10.26448246205318
$ws.Range("M22").Value = 19.02776630132103
$ws.Range("O22").Value = 24.30751665091303
$ws.Range("C23").Value = 10.59552232153538
$ws.Range("D23").Value = 9.282586427418044
$ws.Range("E23").Value = 13.58850647241667
$ws.Range("F23").Value = 32.33558427601773
$ws.Range("G23").Value = 33.55303114056906
$ws.Range("H23").Value = 15.59290076850512
$ws.Range("I23").Value = 25.32274673567316
$ws.Range("J23").Value = 10.10766465550931
$ws.Range("K23").Value = 17.93029830729185
$ws.Range("L23").Value = 10.27311709861122
$ws.Range("M23").Value = 18.92178150522997
$ws.Range("O23").Value = 24.35415226204845
$ws.Range("C24").Value = 10.52514929422026
$ws.Range("D24").Value = 9.227296838319997
$ws.Range("E24").Value = 13.61618123908285
$ws.Range("F24").Value = 32.52949854417091
$ws.Range("G24").Value = 33.84489610753036
$ws.Range("H24").Value = 15.69424536785142
$ws.Range("I24").Value = 25.48166531613598
$ws.Range("J24").Value = 10.14497595924689
$ws.Range("K24").Value = 17.05178294943721
$ws.Range("L24").Value = 10.30730532607158
$ws.Range("M24").Value = 18.51592018377616
$ws.Range("O24").Value = 24.54051496195769
$ws.Range("C25").Value = 10.45429629257827
$ws.Range("D25").Value = 9.170977749671147
$ws.Range("E25").Value = 13.6524248460396
$ws.Range("F25").Value = 32.76860911869425
$ws.Range("G25").Value = 34.20187795468318
$ws.Range("H25").Value = 15.81354207066128
$ws.Range("I25").Value = 25.67414273993565
$ws.Range("J25").Value = 10.18839450663497
$ws.Range("K25").Value = 16.04699525697279
$ws.Range("L25").Value = 10.34737330080899
$ws.Range("M25").Value = 18.07040046121693
$ws.Range("O25").Value = 24.76228120714851
